$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns for changed cells to preserve exact
# string formatting (e.g. trailing zeros, thousand-dot separators) instead of
# letting Excel auto-convert numeric-looking strings into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.168.17"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.893.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3729"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07215"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9068"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.09"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07637"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.890.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.28"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.291"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008519"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9992"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.218.37"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.060"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.138.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.462"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.82"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.794"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.155"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.73"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.914"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.804"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09218"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05066"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.198"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7628"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.031"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.287"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.572"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5633"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02000"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.079"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.82"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.603"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.878"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1510"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4808"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.579"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.16"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.40%  "
